# Applies a permutation of the species-identification / coordinate block
# (columns A, B, D, E, F, G, H, K, L, M, N, Q, R, AC) among data rows 2-20
# of the "Artfynd" sheet. Each row keeps its own observation-event metadata
# (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY); only the
# taxon/coordinate/comment block moves between rows, per the mapping below
# (target row -> source row that its new block comes from).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together as one block.
$blockCols = @("A","B","D","E","F","G","H","Q","R","AC")

# target row -> source row (source row's block becomes target row's new block)
$mapping = @{
    2  = 4
    3  = 15
    4  = 6
    5  = 7
    6  = 17
    7  = 2
    8  = 14
    9  = 11
    10 = 18
    11 = 5
    12 = 12
    13 = 9
    14 = 10
    15 = 20
    16 = 16
    17 = 3
    18 = 19
    19 = 8
    20 = 13
}

# 1) Snapshot every row's current block values BEFORE any writes, so that
#    writes to one row never clobber data still needed as a source for
#    another row.
$snapshot = @{}
foreach ($r in 2..20) {
    $row = @{}
    foreach ($c in $blockCols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    # K/L/M/N are always either all-present(empty) or all-absent together,
    # in lockstep with F/G/H being the "Tretåig hackspett" species block.
    $row["hasKLMN"] = ($ws.Range("F$r").Value2 -eq "Tretåig hackspett")
    $snapshot[$r] = $row
}

# 2) Write each target row's new block from the mapped source row's snapshot.
foreach ($targetRow in 2..20) {
    $sourceRow = $mapping[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    $src = $snapshot[$sourceRow]

    foreach ($c in @("A","B","E","F","G","H","Q","R")) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }

    # D ("Rödlistade") is only populated for some species; mirror source.
    if ($src["D"] -eq $null -or $src["D"] -eq "") {
        $ws.Range("D$targetRow").ClearContents()
    } else {
        $ws.Range("D$targetRow").Value = $src["D"]
    }

    # AC ("Publik kommentar") likewise only populated for some rows.
    if ($src["AC"] -eq $null -or $src["AC"] -eq "") {
        $ws.Range("AC$targetRow").ClearContents()
    } else {
        $ws.Range("AC$targetRow").Value = $src["AC"]
    }

    # K, L, M, N are blank placeholder cells that exist only on
    # "Tretåig hackspett" rows; add/remove them to match the new species.
    if ($src["hasKLMN"]) {
        foreach ($c in @("K","L","M","N")) {
            $ws.Range("$c$targetRow").Value = ""
        }
    } else {
        foreach ($c in @("K","L","M","N")) {
            $ws.Range("$c$targetRow").ClearContents()
        }
    }
}
